# Refresh the cryptos worksheet with the latest price/volume snapshot
# (GitHub Actions scrape update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new text value. Column D entries that look like plain numbers are
# prefixed with a leading apostrophe so Excel stores them as text (matching
# the workbook's original inline-string cells, e.g. keeping "25.90" from
# collapsing to 25.9).
$updates = @(
    @{ Cell = 'D2'; Value = '56.547.65' }
    @{ Cell = 'E2'; Value = '  -2.50%  ' }
    @{ Cell = 'D3'; Value = '2.958.28' }
    @{ Cell = 'E4'; Value = '  +0.04%  ' }
    @{ Cell = 'D5'; Value = '''496.42' }
    @{ Cell = 'E5'; Value = '  -5.93%  ' }
    @{ Cell = 'D6'; Value = '''134.29' }
    @{ Cell = 'E6'; Value = '  -6.48%  ' }
    @{ Cell = 'E7'; Value = '  +0.01%  ' }
    @{ Cell = 'D8'; Value = '''0.424' }
    @{ Cell = 'E8'; Value = '  -5.36%  ' }
    @{ Cell = 'D9'; Value = '''7.14' }
    @{ Cell = 'E9'; Value = '  -6.59%  ' }
    @{ Cell = 'E10'; Value = '  -6.51%  ' }
    @{ Cell = 'E11'; Value = '  -5.37%  ' }
    @{ Cell = 'D12'; Value = '3.465.22' }
    @{ Cell = 'E12'; Value = '  -3.45%  ' }
    @{ Cell = 'E13'; Value = '  -3.29%  ' }
    @{ Cell = 'D14'; Value = '''25.90' }
    @{ Cell = 'E14'; Value = '  -5.51%  ' }
    @{ Cell = 'E15'; Value = '  -9.18%  ' }
    @{ Cell = 'D16'; Value = '56.642.72' }
    @{ Cell = 'E16'; Value = '  -2.31%  ' }
    @{ Cell = 'D17'; Value = '2.957.49' }
    @{ Cell = 'E17'; Value = '  -3.51%  ' }
    @{ Cell = 'D18'; Value = '''5.96' }
    @{ Cell = 'E18'; Value = '  -4.15%  ' }
    @{ Cell = 'D19'; Value = '''12.46' }
    @{ Cell = 'E19'; Value = '  -5.77%  ' }
    @{ Cell = 'D20'; Value = '''7.74' }
    @{ Cell = 'E20'; Value = '  -6.01%  ' }
    @{ Cell = 'D21'; Value = '''316.66' }
    @{ Cell = 'E21'; Value = '  -7.49%  ' }
    @{ Cell = 'D22'; Value = '''0.998' }
    @{ Cell = 'E22'; Value = '  -0.15%  ' }
    @{ Cell = 'D23'; Value = '''5.70' }
    @{ Cell = 'E23'; Value = '  +0.59%  ' }
    @{ Cell = 'D24'; Value = '''0.485' }
    @{ Cell = 'E24'; Value = '  -4.07%  ' }
    @{ Cell = 'D25'; Value = '''62.62' }
    @{ Cell = 'E25'; Value = '  -4.30%  ' }
    @{ Cell = 'E26'; Value = '  +0.22%  ' }
    @{ Cell = 'E27'; Value = '  -5.36%  ' }
    @{ Cell = 'D28'; Value = '0.0₃0861' }
    @{ Cell = 'E28'; Value = '  -12.35%  ' }
    @{ Cell = 'D29'; Value = '''6.47' }
    @{ Cell = 'E29'; Value = '  -7.99%  ' }
    @{ Cell = 'D30'; Value = '''7.00' }
    @{ Cell = 'E30'; Value = '  -6.07%  ' }
    @{ Cell = 'E31'; Value = '  -6.03%  ' }
    @{ Cell = 'D32'; Value = '''19.89' }
    @{ Cell = 'E32'; Value = '  -6.35%  ' }
    @{ Cell = 'E33'; Value = '  -8.81%  ' }
    @{ Cell = 'D34'; Value = '''152.26' }
    @{ Cell = 'E34'; Value = '  -3.49%  ' }
    @{ Cell = 'D35'; Value = '''4.46' }
    @{ Cell = 'E35'; Value = '  -7.57%  ' }
    @{ Cell = 'D36'; Value = '''5.67' }
    @{ Cell = 'E36'; Value = '  -5.38%  ' }
    @{ Cell = 'D37'; Value = '''1.20' }
    @{ Cell = 'E37'; Value = '  -9.80%  ' }
    @{ Cell = 'D38'; Value = '''23.78' }
    @{ Cell = 'E38'; Value = '  -9.36%  ' }
    @{ Cell = 'D39'; Value = '''0.0651' }
    @{ Cell = 'E39'; Value = '  -7.46%  ' }
    @{ Cell = 'D40'; Value = '2.989.05' }
    @{ Cell = 'E40'; Value = '  -3.51%  ' }
    @{ Cell = 'D41'; Value = '''37.26' }
    @{ Cell = 'E41'; Value = '  -1.38%  ' }
    @{ Cell = 'E42'; Value = '  -0.02%  ' }
    @{ Cell = 'E43'; Value = '  -4.17%  ' }
    @{ Cell = 'D44'; Value = '''3.67' }
    @{ Cell = 'E44'; Value = '  -7.48%  ' }
    @{ Cell = 'D45'; Value = '2.149.85' }
    @{ Cell = 'E45'; Value = '  -8.40%  ' }
    @{ Cell = 'E46'; Value = '  -9.01%  ' }
    @{ Cell = 'E47'; Value = '  -4.70%  ' }
    @{ Cell = 'D48'; Value = '''0.922' }
    @{ Cell = 'E48'; Value = '  -12.24%  ' }
    @{ Cell = 'D49'; Value = '''0.0230' }
    @{ Cell = 'E49'; Value = '  -6.15%  ' }
    @{ Cell = 'D50'; Value = '''18.94' }
    @{ Cell = 'E50'; Value = '  -6.90%  ' }
    @{ Cell = 'D51'; Value = '''1.73' }
    @{ Cell = 'E51'; Value = '  -13.93%  ' }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
